# Append the "13.12.2021" diary entry (profile-viewing activity) to the
# end of the learning diary, mirroring how it would be typed in Word.

$d = $word.ActiveDocument

# Start from the very end of the document body (after the last existing
# paragraph, "...pretty much completed and functional at this point.")
$end = $d.Content
$end.Collapse(0)

# --- two blank separator paragraphs -------------------------------------
$end.InsertParagraphAfter()
$end.Collapse(0)
$end.InsertParagraphAfter()
$end.Collapse(0)

# --- date line ------------------------------------------------------------
$end.InsertParagraphAfter()
$end.Collapse(0)
$datePara = $d.Paragraphs.Last.Range
$datePara.Collapse(0)
$datePara.InsertAfter("13.12.2021")

# --- blank line ------------------------------------------------------------
$end2 = $d.Content
$end2.Collapse(0)
$end2.InsertParagraphAfter()
$end2.Collapse(0)

# --- "goal for this session" paragraph -------------------------------------
$end2.InsertParagraphAfter()
$end2.Collapse(0)
$goalPara = $d.Paragraphs.Last.Range
$goalPara.Collapse(0)
$goalPara.InsertAfter("The goal I set for this session was to create the activity for viewing a profile. I started by creating a new empty activity and making a basic layout for the components. Next I started working on the code for fetching profile information. I felt like I could have made the fetcher its own separate class so I could reuse the same code in both activities, but I might do that refactoring later.")

# --- blank line ------------------------------------------------------------
$end3 = $d.Content
$end3.Collapse(0)
$end3.InsertParagraphAfter()
$end3.Collapse(0)

# --- "profile pictures / load image from URL" paragraph (with hyperlink) ---
$end3.InsertParagraphAfter()
$end3.Collapse(0)
$picPara = $d.Paragraphs.Last.Range
$picPara.Collapse(0)
$picPara.InsertAfter("To display profile pictures, I needed to find a way to load an image from an URL. I watched this tutorial to learn how: https://www.youtube.com/watch?v=oz3uGdi3f8Q. After implementing this and trying to test it, I found out that the response body for the API has been recently changed, adding a name to the JSON array. Previously it was unnamed, so I had to change my code to accommodate for this.")

$picParaFull = $d.Paragraphs.Last.Range
$urlRange = $picParaFull.Duplicate
$urlRange.Find.Execute("https://www.youtube.com/watch?v=oz3uGdi3f8Q")
$h1 = $d.Hyperlinks.Add($urlRange, "https://www.youtube.com/watch?v=oz3uGdi3f8Q")
$h1.Range.Style = "Hyperlinkki"

# --- blank line ------------------------------------------------------------
$end4 = $d.Content
$end4.Collapse(0)
$end4.InsertParagraphAfter()
$end4.Collapse(0)

# --- "getText/EditText bug, country flag" paragraph (with hyperlink) -------
$end4.InsertParagraphAfter()
$end4.Collapse(0)
$bugPara = $d.Paragraphs.Last.Range
$bugPara.Collapse(0)
$bugPara.InsertAfter("The next problem was an error that happened because using getText on an EditText doesn")
$bugPara2 = $d.Paragraphs.Last.Range
$bugPara2.Collapse(0)
$bugPara2.InsertAfter([char]0x2019 + "t return a String but an Editable, so I fixed that. After that, I noticed that one of the info pieces I wanted from the JSON was inside another object, so I used getJSONObject twice to get the data. After that I styled the profile activity to be consistent with the main one. I wanted to display a country flag instead of just text on the profiles, so I looked up how to convert a country code to a flag emoji: https://attacomsian.com/blog/how-to-convert-country-code-to-emoji-in-java. ")

$bugParaFull = $d.Paragraphs.Last.Range
$urlRange2 = $bugParaFull.Duplicate
$urlRange2.Find.Execute("https://attacomsian.com/blog/how-to-convert-country-code-to-emoji-in-java")
$h2 = $d.Hyperlinks.Add($urlRange2, "https://attacomsian.com/blog/how-to-convert-country-code-to-emoji-in-java")
$h2.Range.Style = "Hyperlinkki"

Write-Output ("Final paragraph count: " + $d.Paragraphs.Count)
